$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 47, pushing existing rows 47-53 down to 48-54.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new record's data.
$ws.Cells.Item(47, 1).Value = 3
$ws.Cells.Item(47, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44449
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(47, 5).Value = 5
$ws.Cells.Item(47, 6).Value = 100112026
$ws.Cells.Item(47, 7).Value = "Haba"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 55
$ws.Cells.Item(47, 11).Value = 13000
$ws.Cells.Item(47, 12).Value = 14000
$ws.Cells.Item(47, 13).Value = 13545
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(47, 16).Value = 542
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
